$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cookies")

# Header (D7:E7 merged)
$ws.Range("D7").Value2 = "Emissions from cookies - yahoo.com"

# Row 8: Page views per month (B) -- now in billions, value 18.942, with new numeric format; add Source note in G8
$ws.Range("D8").Value2 = "Page views per month (B)"
$ws.Range("E8").Value2 = 18.942
$ws.Range("E8").NumberFormat = "#,##0.0"
$ws.Range("G8").Value2 = "Source: SimilarWeb"

# Row 11: Cookie syncs per month now computed off billions value
$ws.Range("E11").Formula = "=E10*E8"

# Row 13: Data transfer, account for E11 now being expressed in billions
$ws.Range("E13").Formula = "=1.3*E11*1000000000/1024/1024"

# Row 17: Server-side emissions, account for E11 now being expressed in billions
$ws.Range("E17").Formula = "=0.000365*E11*1000000000"

# Row 19: Total annual emissions (multiply by 12 for annual, rename label)
$ws.Range("D19").Value2 = "Total annual emissions (mt CO2e)"
$ws.Range("E19").Formula = "=(E17+E15)/1000000*12"

# Update selection to match target view
$ws.Range("G16").Select() | Out-Null
